# Daily refresh of the cryptos price/volume table (GitHub Actions job).
# Column D ("Price") holds numeric-looking text (e.g. "69.692.57", "1.00")
# that must stay plain text, so those assignments are prefixed with a
# leading apostrophe - the normal Excel "force text" input trick - to stop
# Excel's automatic number coercion from collapsing e.g. "1.00" to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.692.57'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '''3.705.05'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''683.65'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = '''161.08'
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '''0.497'
$ws.Range("E8").Value = '  +0.51%  '
$ws.Range("D9").Value = '''0.147'
$ws.Range("E9").Value = '  +1.10%  '
$ws.Range("D10").Value = '''7.16'
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").Value = '''0.442'
$ws.Range("E11").Value = '  +1.49%  '
$ws.Range("D12").Value = '''0.0000234'
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("D13").Value = '''4.326.12'
$ws.Range("E13").Value = '  +0.56%  '
$ws.Range("D14").Value = '''32.60'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '''3.697.09'
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").Value = '''69.557.87'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("E17").Value = '  +3.09%  '
$ws.Range("D18").Value = '''16.10'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").Value = '''6.48'
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").Value = '''473.87'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("D21").Value = '''9.98'
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("D22").Value = '''0.654'
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("D23").Value = '''80.38'
$ws.Range("E23").Value = '  +1.15%  '
$ws.Range("D24").Value = '''3.846.83'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = '''0.0000126'
$ws.Range("E26").Value = '  -0.75%  '
$ws.Range("D27").Value = '''11.05'
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").Value = '''9.21'
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("D29").Value = '''2.72'
$ws.Range("E29").Value = '  +0.79%  '
$ws.Range("D30").Value = '''1.74'
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").Value = '''2.02'
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = '''6.59'
$ws.Range("E32").Value = '  -1.52%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''27.07'
$ws.Range("E33").Value = '  +1.38%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").Value = '''3.693.06'
$ws.Range("E35").Value = '  +1.18%  '
$ws.Range("D36").Value = '''0.160'
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").Value = '''8.40'
$ws.Range("E37").Value = '  +3.11%  '
$ws.Range("D38").Value = '''6.29'
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '''2.30'
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").Value = '''0.998'
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").Value = '''0.0911'
$ws.Range("E42").Value = '  +0.94%  '
$ws.Range("D43").Value = '''169.36'
$ws.Range("E43").Value = '  +2.63%  '
$ws.Range("D44").Value = '''0.945'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = '''47.36'
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''29.08'
$ws.Range("E46").Value = '  +1.22%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '''2.75'
$ws.Range("E47").Value = '  +1.78%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = '''0.000282'
$ws.Range("E48").Value = '  +3.97%  '
$ws.Range("D49").Value = '''1.12'
$ws.Range("E49").Value = '  +2.53%  '
$ws.Range("D50").Value = '''1.31'
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").Value = '''7.89'
$ws.Range("E51").Value = '  +0.07%  '

Write-Output "Applied cryptos update"
